$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts C:BD -> D:BE), carrying over
# formatting/width of the old column C to the new column D, etc.
$null = $ws.Columns("C").Insert()

# Populate the new "Status" column.
$ws.Range("C1").Value = "Status"
$ws.Range("C2").Value = '${table:rawData.status2}'

# Match the look of the surrounding header/data cells (copy number format,
# font, alignment, etc. from the neighboring column rather than re-creating
# style records from scratch).
$ws.Range("D1").Copy()
$null = $ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$null = $ws.Range("C2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Give the new column its own width.
$ws.Columns("C").ColumnWidth = 22.5

# Reset the view: scroll back to show column A and select D2 (first data cell).
$ws.Activate()
$null = $ws.Range("D2").Select()
